$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '28.646.94'
$ws.Cells.Item(2, 5).Value = "  +4.29%  "
Set-TextValue $ws.Cells.Item(3, 4) '1.871.96'
$ws.Cells.Item(3, 5).Value = "  +2.36%  "
Set-TextValue $ws.Cells.Item(4, 4) '0.9977'
$ws.Cells.Item(4, 5).Value = "  -0.78%  "
Set-TextValue $ws.Cells.Item(5, 4) '338.54'
$ws.Cells.Item(5, 5).Value = "  +2.05%  "
Set-TextValue $ws.Cells.Item(6, 4) '0.9985'
$ws.Cells.Item(6, 5).Value = "  -0.70%  "
Set-TextValue $ws.Cells.Item(7, 4) '0.4698'
$ws.Cells.Item(7, 5).Value = "  +2.97%  "
Set-TextValue $ws.Cells.Item(8, 4) '0.3987'
$ws.Cells.Item(8, 5).Value = "  +4.48%  "
Set-TextValue $ws.Cells.Item(9, 4) '47.69'
$ws.Cells.Item(9, 5).Value = "  +3.02%  "
Set-TextValue $ws.Cells.Item(10, 4) '0.08056'
$ws.Cells.Item(10, 5).Value = "  +1.92%  "
Set-TextValue $ws.Cells.Item(11, 4) '1.003'
$ws.Cells.Item(11, 5).Value = "  +3.33%  "
Set-TextValue $ws.Cells.Item(12, 4) '22.14'
$ws.Cells.Item(12, 5).Value = "  +5.44%  "
Set-TextValue $ws.Cells.Item(13, 4) '6.060'
$ws.Cells.Item(13, 5).Value = "  +3.25%  "
Set-TextValue $ws.Cells.Item(14, 4) '1.863.69'
$ws.Cells.Item(14, 5).Value = "  +1.93%  "
Set-TextValue $ws.Cells.Item(15, 4) '7.307'
$ws.Cells.Item(15, 5).Value = "  +3.71%  "
Set-TextValue $ws.Cells.Item(16, 4) '90.47'
$ws.Cells.Item(16, 5).Value = "  +1.93%  "
Set-TextValue $ws.Cells.Item(17, 4) '1.001'
$ws.Cells.Item(17, 5).Value = "  -0.46%  "
Set-TextValue $ws.Cells.Item(18, 4) '0.00001040'
$ws.Cells.Item(18, 5).Value = "  +1.00%  "
Set-TextValue $ws.Cells.Item(19, 4) '0.06627'
$ws.Cells.Item(19, 5).Value = "  -0.23%  "
Set-TextValue $ws.Cells.Item(20, 4) '17.59'
$ws.Cells.Item(20, 5).Value = "  +2.27%  "
Set-TextValue $ws.Cells.Item(21, 4) '0.9999'
$ws.Cells.Item(21, 5).Value = "  -0.52%  "
Set-TextValue $ws.Cells.Item(22, 4) '28.587.21'
$ws.Cells.Item(22, 5).Value = "  +4.22%  "
Set-TextValue $ws.Cells.Item(23, 4) '5.489'
$ws.Cells.Item(23, 5).Value = "  +3.03%  "
Set-TextValue $ws.Cells.Item(24, 4) '11.05'
$ws.Cells.Item(24, 5).Value = "  +2.41%  "
Set-TextValue $ws.Cells.Item(25, 4) '2.255'
$ws.Cells.Item(25, 5).Value = "  -2.24%  "
Set-TextValue $ws.Cells.Item(26, 4) '2.082.16'
$ws.Cells.Item(26, 5).Value = "  +1.62%  "
Set-TextValue $ws.Cells.Item(27, 4) '161.12'
$ws.Cells.Item(27, 5).Value = "  +2.60%  "
Set-TextValue $ws.Cells.Item(28, 4) '19.79'
$ws.Cells.Item(28, 5).Value = "  +2.06%  "
Set-TextValue $ws.Cells.Item(29, 4) '2.118'
$ws.Cells.Item(29, 5).Value = "  +2.66%  "
Set-TextValue $ws.Cells.Item(30, 4) '5.495'
$ws.Cells.Item(30, 5).Value = "  +4.88%  "
Set-TextValue $ws.Cells.Item(31, 4) '119.91'
$ws.Cells.Item(31, 5).Value = "  +1.50%  "
Set-TextValue $ws.Cells.Item(32, 4) '0.9725'
$ws.Cells.Item(32, 5).Value = "  +2.70%  "
Set-TextValue $ws.Cells.Item(33, 4) '0.09539'
$ws.Cells.Item(33, 5).Value = "  +2.57%  "
Set-TextValue $ws.Cells.Item(34, 4) '3.589'
$ws.Cells.Item(34, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(37, 4) '0.06232'
$ws.Cells.Item(37, 5).Value = "  +5.07%  "
Set-TextValue $ws.Cells.Item(38, 4) '0.02253'
$ws.Cells.Item(38, 5).Value = "  +3.55%  "
Set-TextValue $ws.Cells.Item(39, 4) '8.371'
$ws.Cells.Item(39, 5).Value = "  +4.36%  "
Set-TextValue $ws.Cells.Item(40, 4) '1.183'
$ws.Cells.Item(40, 5).Value = "  +2.02%  "
Set-TextValue $ws.Cells.Item(41, 4) '0.5936'
$ws.Cells.Item(41, 5).Value = "  +3.14%  "
Set-TextValue $ws.Cells.Item(42, 4) '0.9993'
$ws.Cells.Item(42, 5).Value = "  -0.65%  "
Set-TextValue $ws.Cells.Item(43, 4) '0.1882'
$ws.Cells.Item(43, 5).Value = "  +2.69%  "
Set-TextValue $ws.Cells.Item(44, 4) '10.35'
$ws.Cells.Item(44, 5).Value = "  +3.34%  "
Set-TextValue $ws.Cells.Item(45, 4) '1.259'
$ws.Cells.Item(45, 5).Value = "  -0.14%  "
Set-TextValue $ws.Cells.Item(49, 4) '1.955'
$ws.Cells.Item(49, 5).Value = "  +4.73%  "
Set-TextValue $ws.Cells.Item(50, 4) '2.080'
$ws.Cells.Item(50, 5).Value = "  +13.41%  "
Set-TextValue $ws.Cells.Item(51, 4) '112.52'
$ws.Cells.Item(51, 5).Value = "  +2.05%  "

# Rows with coin re-ordering (B, C, D, E all change)
$ws.Cells.Item(35, 2).Value = "ARBITRUM"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(35, 4) '1.386'
$ws.Cells.Item(35, 5).Value = "  +5.00%  "
$ws.Cells.Item(36, 2).Value = "Filecoin"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(36, 4) '5.374'
$ws.Cells.Item(36, 5).Value = "  +2.58%  "
$ws.Cells.Item(46, 2).Value = "Decentraland"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Cells.Item(46, 4) '0.5572'
$ws.Cells.Item(46, 5).Value = "  +2.14%  "
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(47, 4) '12.16'
$ws.Cells.Item(47, 5).Value = "  +1.21%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(48, 4) '0.07415'
$ws.Cells.Item(48, 5).Value = "  +11.95%  "
